$d = $word.ActiveDocument

# --- document.xml: fill the two trailing empty paragraphs with text,
# then append two more empty paragraphs (one plain, one carrying an
# explicit copy of style0's paragraph formatting). ---

$p6 = $d.Paragraphs(6)
$p6.Range.Text = "links!"

$p7 = $d.Paragraphs(7)
$p7.Range.Text = "Holo theme switching"

# New, plain empty paragraph.
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)

# Another new, empty paragraph that gets explicit direct formatting
# mirroring style0 (jc/widowControl/tabs/suppressAutoHyphens/spacing).
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs(9)

# Alignment: toggle away and back via the paragraph-collection so the
# (otherwise redundant-with-style) "left" value is written explicitly.
$paras = $p9.Range.Paragraphs
$paras.Alignment = 1
$paras.Alignment = 0

$p9.Format.WidowControl = 1

$ts = $p9.Format.TabStops.Add(35.45)
$ts.Leader = 0

$p9.Format.Hyphenation = $false

$p9.Format.SpaceAfter = 10
$p9.Format.SpaceBefore = 0
$p9.Format.LineSpacingRule = 3
$p9.Format.LineSpacing = 13.8

# --- styles.xml: style0's default run color auto -> 00000A ---
$st = $d.Styles("style0")
$st.Font.Color = 655360
